$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("G2").Value = 1.5
$ws.Range("H2").Value = 3.9
$ws.Range("J2").Value = 2.1
$ws.Range("K2").Value = 2.1
$ws.Range("L2").Value = 7
$ws.Range("M2").Value = 1.08
$ws.Range("N2").Value = 8
$ws.Range("Z2").Value = 10
$ws.Range("AC2").Value = 8
$ws.Range("AD2").Value = 7.5
$ws.Range("AG2").Value = 13
$ws.Range("AI2").Value = 21
$ws.Range("AM2").Value = 3.25
$ws.Range("AP2").Value = 26
$ws.Range("AY2").Value = 151

# Row 3 updates
$ws.Range("BC3").Value = 151

# Row 4 updates
$ws.Range("M4").Value = 1.14
$ws.Range("N4").Value = 5.5
